$d = $word.ActiveDocument

# --- Fix the "staff name" introduction line -------------------------------
# The three runs "_______________, [staff name]" + "," + " " (the latter two
# wrapped by proofErr gramStart/gramEnd markers) become a single italic run
# "_______________, [staff name], ".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "_______________, [staff name], "
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "_______________, [staff name], "
$find.Execute([ref]"_______________, [staff name], ", $false, $false, $false, $false, $false, $true, $true, $false, "_______________, [staff name], ", 2)

# --- Finalize tracked changes ---------------------------------------------
# Accept the pending insertion ("study ") and deletions (the old
# "Treatment and compensation for injury" paragraphs), which also finalizes
# the run-merge replacement above, leaving no revision marks behind.
$d.Revisions.AcceptAll()

# --- Reposition the _GoBack bookmark ---------------------------------------
# Redefine _GoBack so it spans from the very start of the document to just
# before the final empty paragraph (i.e. right after the signature block),
# matching where Word leaves it after the last edit in the document.
$endPos = $d.Paragraphs($d.Paragraphs.Count - 1).Range.Start
$goBackRange = $d.Range(0, $endPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
